$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

$ws.Range("D2").Value = "No"
$ws.Range("D3").Value = "No"
$ws.Range("D4").Value = "No"
$ws.Range("D5").Value = "Yes"
$ws.Range("D6").Value = "Yes"

$ws.Range("D6").Select()
